$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistic")

# Update existing row 1 values (A1 and D1 are unchanged)
$ws.Cells.Item(1, 2).Value = 0
$ws.Cells.Item(1, 3).Value = 7.9149503000000001

# New data rows 2-8 (columns A, B, C only)
$data = @(
    @(2, 0, 7.4623390000000001),
    @(3, 0, 7.3745539999999998),
    @(4, 0, 7.876023),
    @(5, 0, 7.5219833999999999),
    @(6, 0, 8.0597259000000001),
    @(7, 0, 7.4374422999999998),
    @(8, 0, 8.0246790000000008)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Adjust column C width to match column D (9.7109375)
$ws.Columns.Item(3).ColumnWidth = 8.8
